$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

# Columns A (Vaccine) and D (Packaging) contain stray HTML-entity artifacts
# (&curren; / &bull;), a stray "#" character, and embedded line breaks that
# should be collapsed to single spaces.
$columns = @(1, 4)

for ($r = 1; $r -le $lastRow; $r++) {
    foreach ($col in $columns) {
        $cell = $ws.Cells.Item($r, $col)
        $val = $cell.Text
        if ($val) {
            $newVal = $val.Replace("&curren;", "").Replace("&bull;", "").Replace("#", "").Replace("`n", " ")
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
